$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: survival/error formulas switched from a simple multiply to the
# (count / (sightability - SD)) - paired-estimate form.
$ws.Range("L3").Formula = "=(F3/(I3-J3))-K3"
$ws.Range("N3").Formula = "=(G3/(I3-J3))-M3"

# Row 9
$ws.Range("L9").Formula = "=(F9/(I9-J9))-K9"
$ws.Range("N9").Formula = "=(G9/(I9-J9))-M9"

# Row 14
$ws.Range("L14").Formula = "=(F14/(I14-J14))-K14"
$ws.Range("N14").Formula = "=(G14/(I14-J14))-M14"

# Row 17 - only the N column error formula changed (L17 left as-is)
$ws.Range("N17").Formula = "=(G17/(I17-J17))-M17"

# Row 19 - I19 gains an explicit formula (still evaluates to 0.75)
$ws.Range("I19").Formula = "=0.75"
$ws.Range("L19").Formula = "=(F19/(I19-J19))-K19"
$ws.Range("N19").Formula = "=(G19/(I19-J19))-M19"

# Row 20
$ws.Range("L20").Formula = "=(F20/(I20-J20))-K20"
$ws.Range("N20").Formula = "=(G20/(I20-J20))-M20"

# Row 22
$ws.Range("L22").Formula = "=(F22/(I22-J22))-K22"
$ws.Range("N22").Formula = "=(G22/(I22-J22))-M22"

# Scroll/selection state: the sheet view now shows row 3 at the top with
# N29 selected (was M22 with no scroll offset).
$ws.Range("N29").Select()
